$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.36509669890178031,
    -0.0099999995999482394,
    -0.0089999995942591227,
    0.061996261466966018,
    -0.0059999995991528721,
    -0.020927391766896619,
    -0.019999999507115618,
    -0.019999999503421684,
    -0.0059999995743531542,
    -0.020151785404735278,
    -0.0044999995794299252,
    -0.0059999995713737597,
    -0.0059999995717570087,
    -0.011999999540277528,
    -0.0059999995732251676,
    -0.0059999995750392721,
    -0.005999999579210602,
    0.039182502775344297,
    -0.084191257677313391,
    -0.0089999996000358351,
    -0.0089999995994523019,
    -0.0089999995990508452,
    -0.0089999995906353547,
    -0.041999999408319333,
    -0.041999999404983335,
    -0.0059999995811104156,
    -0.00599999957865327,
    -0.0059999995695809716,
    -0.011999999531139949,
    -0.019999999485459607,
    -0.014999999507582018,
    -0.020999999474978992,
    -0.0059999995545672036
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
